# Update existing data on Sheet1: new values for rows 2-13 (cols A,B,E,F),
# delete rows 14-23, then add a new empty Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @{ A = -1.292080721726155;  B = 40.47579836825359; E = -29.97; F = 19.93 }
    3  = @{ A = -32.37101184192168;  B = 18.50804372641845; E = -30.1;  F = 19.53 }
    4  = @{ A = -32.354482725765;    B = 18.735814354647;   E = -29.95; F = 20.03 }
    5  = @{ A = -32.3063168398515;   B = 18.98166524023402; E = -30.06; F = 19.58 }
    6  = @{ A = -32.35458014444087;  B = 18.94176099632543; E = -30.04; F = 19.76 }
    7  = @{ A = -32.31753401087381;  B = 18.8663606345545;  E = -30;    F = 19.88 }
    8  = @{ A = -32.3397456370814;   B = 18.99823253596354; E = -30.07; F = 19.57 }
    9  = @{ A = -32.3719895303408;   B = 18.57676856475656; E = -29.96; F = 20.16 }
    10 = @{ A = -32.28325829338218;  B = 19.41418914595164; E = -30.09; F = 19.75 }
    11 = @{ A = -32.26272399813831;  B = 19.45634640586082; E = -30.15; F = 19.66 }
    12 = @{ A = -32.27723217532559;  B = 19.50693954850279; E = -30.19; F = 19.6  }
    13 = @{ A = -32.27415980015615;  B = 19.35380986947552 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    if ($row.ContainsKey("A")) { $ws.Cells.Item($r, 1).Value = $row.A }
    if ($row.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $row.B }
    if ($row.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $row.E }
    if ($row.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $row.F }
}

# Clear E13:F13 (row 13 no longer has these values) and rows 14-23 entirely.
$ws.Range("E13:F13").ClearContents()
$ws.Range("A14:F23").ClearContents()

# Add the new empty Sheet2 right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet2"

# Keep Sheet1 as the active/selected sheet (matches original tabSelected state).
$ws.Activate()
